$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in C4:C6 while preserving their existing style/formatting
$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()

# Update the active selection to D6
$ws.Range("D6").Select()
